$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Transaksi Agen")

# Fix Export Data Transaksi: set the Bank value for row 2 (was blank)
$ws.Range("F2").Value = "BCA"
